$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = "MCH196-1"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 24A | GRAP COUNT NUMER: NONE"
$ws.Range("A2,C2,D2,E2,F2,G2,H2").Font.Name = "Calibri"
$ws.Range("A2,C2,D2,E2,F2,G2,H2").Font.Size = 10
